$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Total horas" summary cells (formula text first, then the total, to
#     match the shared-string insertion order of the original edit) ---
$ws.Range("A24").Value = "4+3+1.5+4+4+3+2+3.5+4.5+4.5+5+1+5+3+4+5+4+4+3+4+3.5+8+3.5+5+5+5+4"
$ws.Range("A1").Value = "Total horas: 106"

# --- Extend the date header row (row 2) with three more days ---
$ws.Range("Z2").Copy()
$ws.Range("AB2:AD2").PasteSpecial(-4122)
$ws.Range("AB2").Value = 44020
$ws.Range("AC2").Value = 44021
$ws.Range("AD2").Value = 44023

# --- Row 4: "Estructura y lectura/escritura XML" gains a 1.5 h. mark in AD ---
$ws.Range("B4").Copy()
$ws.Range("AD4").PasteSpecial(-4122)
$ws.Range("AD4").Value = "1.5 h."

# --- Row 13: "Algoritmo Q-learning" gains 1 h. and 1.5 h. marks in AC/AD ---
$ws.Range("K13").Copy()
$ws.Range("AC13").PasteSpecial(-4122)
$ws.Range("AC13").Value = "1 h."

$ws.Range("L13").Copy()
$ws.Range("AD13").PasteSpecial(-4122)
$ws.Range("AD13").Value = "1.5 h."

# --- Row 17: "Servidor" gains a 3 h. mark in AC ---
$ws.Range("Z17").Copy()
$ws.Range("AC17").PasteSpecial(-4122)
$ws.Range("AC17").Value = "3 h."

# --- Row 20: "Implementación con JPA e Hibernate" gains a 5 h. mark in AB ---
$ws.Range("AA20").Copy()
$ws.Range("AB20").PasteSpecial(-4122)
$ws.Range("AB20").Value = "5 h."

# --- Update sheet view (scrolled/selected cell) ---
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("R27").Select()
